$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $c = $ws.Range($cellAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "27.967.52"
Set-TextValue "E2" "  -3.78%  "
Set-TextValue "D3" "1.873.11"
Set-TextValue "E3" "  -2.86%  "
Set-TextValue "D4" "1.003"
Set-TextValue "E4" "  -0.24%  "
Set-TextValue "D5" "318.20"
Set-TextValue "E5" "  -2.34%  "
Set-TextValue "E6" "  -0.22%  "
Set-TextValue "D7" "0.4361"
Set-TextValue "E7" "  -5.34%  "
Set-TextValue "D8" "0.3757"
Set-TextValue "E8" "  -1.93%  "
Set-TextValue "D9" "0.07473"
Set-TextValue "E9" "  -3.59%  "
Set-TextValue "D10" "0.9357"
Set-TextValue "E10" "  -4.55%  "
Set-TextValue "D11" "21.31"
Set-TextValue "E11" "  -5.62%  "
Set-TextValue "D12" "1.837.64"
Set-TextValue "E12" "  -5.26%  "
Set-TextValue "D13" "6.751"
Set-TextValue "E13" "  -3.30%  "
Set-TextValue "D14" "5.443"
Set-TextValue "E14" "  -4.41%  "
Set-TextValue "D15" "0.06870"
Set-TextValue "E15" "  -2.26%  "
Set-TextValue "D16" "1.003"
Set-TextValue "E16" "  -0.30%  "
Set-TextValue "D17" "81.51"
Set-TextValue "E17" "  -3.47%  "
Set-TextValue "D18" "0.000009049"
Set-TextValue "E18" "  -5.37%  "
Set-TextValue "D19" "1.001"
Set-TextValue "E19" "  -0.25%  "
Set-TextValue "D20" "15.86"
Set-TextValue "E20" "  -5.48%  "
Set-TextValue "D21" "27.975.09"
Set-TextValue "E21" "  -3.89%  "
Set-TextValue "D22" "5.129"
Set-TextValue "D23" "11.06"
Set-TextValue "E23" "  +0.78%  "
Set-TextValue "D24" "2.103.50"
Set-TextValue "E24" "  -3.32%  "
Set-TextValue "D25" "2.036"
Set-TextValue "E25" "  -2.02%  "
Set-TextValue "E26" "  -3.26%  "
Set-TextValue "D27" "18.61"
Set-TextValue "E27" "  -2.41%  "
Set-TextValue "D28" "5.627"
Set-TextValue "E28" "  -0.96%  "
Set-TextValue "D29" "113.43"
Set-TextValue "E29" "  -3.76%  "
Set-TextValue "D30" "1.703"
Set-TextValue "E30" "  -7.81%  "
Set-TextValue "D31" "0.09019"
Set-TextValue "E31" "  -3.28%  "
Set-TextValue "D32" "0.8121"
Set-TextValue "E32" "  -6.17%  "
Set-TextValue "D33" "4.812"
Set-TextValue "E33" "  -6.25%  "
Set-TextValue "D34" "1.182"
Set-TextValue "E34" "  -5.40%  "
Set-TextValue "D35" "2.970"
Set-TextValue "E35" "  -1.69%  "
Set-TextValue "D36" "1.001"
Set-TextValue "E36" "  -0.28%  "
Set-TextValue "D37" "0.05513"
Set-TextValue "E37" "  -3.47%  "
Set-TextValue "D38" "1.117"
Set-TextValue "E38" "  -3.63%  "
Set-TextValue "D39" "0.01981"
Set-TextValue "E39" "  -3.39%  "
Set-TextValue "D40" "2.970"
Set-TextValue "E40" "  -2.83%  "
Set-TextValue "D41" "0.5275"
Set-TextValue "E41" "  -4.55%  "
Set-TextValue "D42" "0.1702"
Set-TextValue "E42" "  -3.02%  "
Set-TextValue "D43" "7.006"
Set-TextValue "E43" "  -7.17%  "
Set-TextValue "D44" "8.777"
Set-TextValue "E44" "  -6.38%  "
Set-TextValue "D45" "0.06760"
Set-TextValue "E45" "  -2.47%  "
Set-TextValue "D46" "0.4892"
Set-TextValue "E46" "  -6.10%  "
Set-TextValue "D47" "10.60"
Set-TextValue "E47" "  -5.67%  "
Set-TextValue "D48" "106.82"
Set-TextValue "E48" "  -3.31%  "
Set-TextValue "D49" "1.676"
Set-TextValue "E49" "  -5.95%  "
Set-TextValue "D50" "1.911"
Set-TextValue "E50" "  -13.66%  "
Set-TextValue "D51" "1.000"
Set-TextValue "E51" "  -0.31%  "
